$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching the existing header style
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Row data: row index, I value, J value
$data = @(
  @(2, 8, 8),
  @(3, 8, 9),
  @(4, 9, 9),
  @(5, 9, 9),
  @(6, 9, 9),
  @(7, 9, 9),
  @(8, 9, 9),
  @(9, 8, 9),
  @(10, 9, 9),
  @(11, 9, 9),
  @(12, 8, 9),
  @(13, 9, 9),
  @(14, 9, 9),
  @(15, 9, 9),
  @(16, 10, 11),
  @(17, 8, 9),
  @(18, 9, 9),
  @(19, 8, 8),
  @(20, 8, 8),
  @(21, 8, 8),
  @(22, 7, 8),
  @(23, 8, 8),
  @(24, 8, 8),
  @(25, 8, 8),
  @(26, 12, 12),
  @(27, 7, 7),
  @(28, 8, 8),
  @(29, 7, 8),
  @(30, 8, 8),
  @(31, 9, 9),
  @(32, 7, 8),
  @(33, 8, 8),
  @(34, 8, 8),
  @(35, 7, 8),
  @(36, 8, 8),
  @(37, 8, 8),
  @(38, 8, 8),
  @(39, 8, 8),
  @(40, 8, 9),
  @(41, 8, 8),
  @(42, 8, 8),
  @(43, 8, 8),
  @(44, 8, 8),
  @(45, 8, 8),
  @(46, 6, 6),
  @(47, 7, 7),
  @(48, 8, 8),
  @(49, 8, 8),
  @(50, 9, 9),
  @(51, 7, 8),
  @(52, 7, 8),
  @(53, 7, 7),
  @(54, 7, 7),
  @(55, 9, 9),
  @(56, 8, 8),
  @(57, 8, 8),
  @(58, 7, 7),
  @(59, 7, 7),
  @(60, 7, 7),
  @(61, 6, 7),
  @(62, 12, 12),
  @(63, 6, 6),
  @(64, 8, 8),
  @(65, 8, 8),
  @(66, 9, 9),
  @(67, 8, 8),
  @(68, 6, 7),
  @(69, 9, 10),
  @(70, 7, 7),
  @(71, 9, 9),
  @(72, 6, 6),
  @(73, 5, 6),
  @(74, 7, 7),
  @(75, 8, 8),
  @(76, 7, 7),
  @(77, 7, 7),
  @(78, 8, 8),
  @(79, 7, 7),
  @(80, 6, 6),
  @(81, 6, 6),
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
